$d = $word.ActiveDocument

# Locate the paragraph that contains the old RFAaffidavit reference inside the
# "{{ field.overflow_value(overflow_message= ... ) }}" Jinja expression.
$paras = $d.Paragraphs
$count = $paras.Count
$target = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*RFAaffidavit.default_overflow_message*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found"
} else {
    $full = $target.Range
    # Exclude the trailing paragraph-mark character so we only replace the
    # visible run content, not the paragraph itself.
    $body = $d.Range($full.Start, $full.End - 1)

    $rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

    $run1 = '<w:r>' + $rPr + '<w:t>{{ field.overflow_value(overflow_message=</w:t></w:r>'
    $run2 = '<w:r>' + $rPr + '<w:t>divorce_no_kids_attachment</w:t></w:r>'
    $run3 = '<w:r>' + $rPr + '<w:t>.</w:t></w:r>'
    $run4 = '<w:r>' + $rPr + '<w:t>default_overflow_message) }}</w:t></w:r>'

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $run1 + $run2 + $run3 + $run4 + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $body.InsertXML($xml)

    Write-Host "Updated: $($target.Range.Text)"
}
